# Updates cryptos list values (prices and 1h volume %) per upstream data refresh.
# Two coin pairs (rows 21/22 and 32/33) were also reordered; their Coin/Link/Price/Volume
# cells are rewritten in place to reflect the swapped ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.609.80"
$ws.Range("E2").Value = "  +3.51%  "
$ws.Range("D3").Value = "3.462.12"
$ws.Range("E3").Value = "  +4.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'578.92"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").Value = "'157.01"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.466.90"
$ws.Range("E8").Value = "  +3.89%  "
$ws.Range("D9").Value = "'0.557"
$ws.Range("E9").Value = "  +4.91%  "
$ws.Range("D10").Value = "'7.60"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").Value = "  +6.06%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("D13").Value = "4.058.47"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'0.0000198"
$ws.Range("E15").Value = "  +9.09%  "
$ws.Range("D16").Value = "'27.81"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "64.624.49"
$ws.Range("E17").Value = "  +3.66%  "
$ws.Range("D18").Value = "3.464.49"
$ws.Range("E18").Value = "  +4.95%  "
$ws.Range("D19").Value = "'6.45"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'14.41"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'398.76"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'8.59"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "'0.547"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").Value = "'72.94"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'0.0000121"
$ws.Range("E26").Value = "  +24.66%  "
$ws.Range("D27").Value = "'9.47"
$ws.Range("E27").Value = "  +6.20%  "
$ws.Range("D28").Value = "'0.181"
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "'5.97"
$ws.Range("E30").Value = "  +7.75%  "
$ws.Range("D31").Value = "'2.05"
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.71"
$ws.Range("E32").Value = "  +5.87%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  +6.62%  "
$ws.Range("D34").Value = "'23.83"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'7.00"
$ws.Range("E36").Value = "  +3.69%  "
$ws.Range("D37").Value = "'1.50"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "'160.89"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "'0.0784"
$ws.Range("E39").Value = "  +7.23%  "
$ws.Range("D40").Value = "'27.70"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").Value = "'1.88"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "2.909.24"
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").Value = "'0.0324"
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "'0.777"
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").Value = "'4.42"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").Value = "'41.81"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("D47").Value = "'1.10"
$ws.Range("E47").Value = "  +5.69%  "
$ws.Range("D48").Value = "'23.00"
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("E49").Value = "  +22.67%  "
$ws.Range("D50").Value = "'0.861"
$ws.Range("E50").Value = "  +5.45%  "
$ws.Range("D51").Value = "'6.56"
$ws.Range("E51").Value = "  +3.77%  "
